$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B. This shifts the existing
# Total_Weight / Population / PerCapWeight columns (B:D) to E:G,
# carrying over their values, shared-string references and styles.
$ws.Range("B1:D1").EntireColumn.Insert()

# The inserted header cells (B1:D1) need the same bold/border style as
# the other header cells (copy format only, not values, from E1).
$ws.Range("E1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The inserted data cells (B2:D9) picked up formatting from the insert;
# the target file has no explicit style on these cells, so clear it.
$ws.Range("B2:D9").ClearFormats()

# New header labels
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "FIPS"
$ws.Range("D1").Value = "Year"

# Row 2 (count)
$ws.Range("B2").Value = 3548
$ws.Range("C2").Value = 3548
$ws.Range("D2").Value = 3548
$ws.Range("E2").Value = 3548
$ws.Range("F2").Value = 3548
$ws.Range("G2").Value = 3548

# Row 3 (mean)
$ws.Range("B3").Value = 1773.5
$ws.Range("C3").Value = 23018.95180383314
$ws.Range("D3").Value = 2009.000563697858
$ws.Range("E3").Value = 48278.40100452794
$ws.Range("F3").Value = 92301.42390078917
$ws.Range("G3").Value = 41952.45182851247

# Row 4 (std)
$ws.Range("B4").Value = 1024.363704940779
$ws.Range("C4").Value = 16271.50384454539
$ws.Range("D4").Value = 2.00112731638899
$ws.Range("E4").Value = 156986.4201255817
$ws.Range("F4").Value = 204117.6904437121
$ws.Range("G4").Value = 31513.92345649195

# Row 5 (min)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1001
$ws.Range("D5").Value = 2006
$ws.Range("E5").Value = 0.3027
$ws.Range("F5").Value = 2893
$ws.Range("G5").Value = 10.46318700311096

# Row 6 (25%)
$ws.Range("B6").Value = 886.75
$ws.Range("C6").Value = 12121
$ws.Range("D6").Value = 2007
$ws.Range("E6").Value = 4341.690424250019
$ws.Range("F6").Value = 16846
$ws.Range("G6").Value = 21230.24606345244

# Row 7 (50%)
$ws.Range("B7").Value = 1773.5
$ws.Range("C7").Value = 13253
$ws.Range("D7").Value = 2009
$ws.Range("E7").Value = 11998.83888182528
$ws.Range("F7").Value = 29479.5
$ws.Range("G7").Value = 34115.93639706657

# Row 8 (75%)
$ws.Range("B8").Value = 2660.25
$ws.Range("C8").Value = 45027
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 35026.13619929253
$ws.Range("F8").Value = 74864.75
$ws.Range("G8").Value = 52642.51799479323

# Row 9 (max)
$ws.Range("B9").Value = 3547
$ws.Range("C9").Value = 47189
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 3026736.844583195
$ws.Range("F9").Value = 2576554
$ws.Range("G9").Value = 263276.1895904924
